# Update "datos actualizados" timestamp in A1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 24 de Julio de 2020 a las 11:11"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 4170328
$ws.Range("C4").Value = 337
$ws.Range("D4").Value = 1980432
$ws.Range("E4").Value = 2042555
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 8
$ws.Range("H4").Value = 147341

# --- Row 20: Banglades ---
$ws.Range("B20").Value = 218658
$ws.Range("C20").Value = 2548
$ws.Range("D20").Value = 120976
$ws.Range("E20").Value = 94846
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 35
$ws.Range("H20").Value = 2836

# --- Row 33: Filipinas ---
$ws.Range("B33").Value = 76444
$ws.Range("C33").Value = 2103
$ws.Range("D33").Value = 24502
$ws.Range("E33").Value = 50063
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 15
$ws.Range("H33").Value = 1879

# --- Rows 40/41: Israel overtakes Emiratos Arabes Unidos in ranking ---
# Row 40 becomes Israel's updated data; row 41 becomes the old (unchanged)
# Emiratos Arabes Unidos data, pushed one row down.
$ws.Range("A40").Value = "Israel"
$ws.Range("B40").Value = 58559
$ws.Range("C40").Value = 577
$ws.Range("D40").Value = 25887
$ws.Range("E40").Value = 32226
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 4
$ws.Range("H40").Value = 446

$ws.Range("A41").Value = "Emiratos Arabes Unidos"
$ws.Range("B41").Value = 57988
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 50848
$ws.Range("E41").Value = 6798
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 342

# --- Row 48: Polonia ---
$ws.Range("B48").Value = 42038
$ws.Range("C48").Value = 458
$ws.Range("D48").Value = 31997
$ws.Range("E48").Value = 8386
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 4
$ws.Range("H48").Value = 1655

# --- Row 51: Barein ---
$ws.Range("E51").Value = 3448
$ws.Range("G51").Value = 2
$ws.Range("H51").Value = 136

# --- Row 64: Austria ---
$ws.Range("B64").Value = 20214
$ws.Range("C64").Value = 115
$ws.Range("D64").Value = 18042
$ws.Range("E64").Value = 1461

# --- Row 123: Eslovaquia ---
$ws.Range("B123").Value = 2118
$ws.Range("C123").Value = 29
$ws.Range("D123").Value = 1577
$ws.Range("E123").Value = 513

# --- Row 124: Eslovenia ---
$ws.Range("B124").Value = 2052
$ws.Range("C124").Value = 19
$ws.Range("D124").Value = 1678
$ws.Range("E124").Value = 259

# --- Row 206: Laos ---
$ws.Range("B206").Value = 20
$ws.Range("C206").Value = 1
$ws.Range("E206").Value = 1

# --- Rows 210/211: Groenlandia / Islas Malvinas swap (tie at 13 cases, order flips) ---
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("A211").Value = "Groenlandia"
